# Rename the "Stock" column header to "Holding" (commit: "changed stock to
# holding"). The header row reads Date, Price, Quantity, Cost, Stock, Total -
# we find the cell by its text rather than a hard-coded address so the edit
# is robust even if columns get reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Text -eq "Stock") {
        $cell.Value = "Holding"
    }
}
